$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '52.088.53'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.990.09'
$ws.Range('E3').Value = '  +1.15%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '354.08'
$ws.Range('E5').Value = '  +0.23%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '107.54'
$ws.Range('E6').Value = '  -4.28%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  +0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.626'
$ws.Range('E9').Value = '  -1.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.22'
$ws.Range('E10').Value = '  -3.32%  '

$ws.Range('E11').Value = '  +1.64%  '

$ws.Range('E12').Value = '  -3.82%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.34'
$ws.Range('E13').Value = '  -3.40%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.462.25'
$ws.Range('E14').Value = '  +1.49%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.69'
$ws.Range('E15').Value = '  -2.26%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.998.01'
$ws.Range('E16').Value = '  +2.15%  '

$ws.Range('E17').Value = '  +3.63%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '52.180.15'
$ws.Range('E18').Value = '  +0.53%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.48'
$ws.Range('E19').Value = '  +4.35%  '

$ws.Range('E20').Value = '  -1.59%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.63'
$ws.Range('E21').Value = '  -5.60%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0971'
$ws.Range('E22').Value = '  -1.70%  '

$ws.Range('E23').Value = '  -2.62%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.99'
$ws.Range('E24').Value = '  -2.03%  '

$ws.Range('E25').Value = '  -2.61%  '

$ws.Range('E26').Value = '  -0.22%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.88'
$ws.Range('E27').Value = '  -1.15%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.58'
$ws.Range('E28').Value = '  +2.15%  '

$ws.Range('E29').Value = '  -0.08%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.109'
$ws.Range('E30').Value = '  -4.47%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.30'
$ws.Range('E31').Value = '  -3.64%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.21'
$ws.Range('E32').Value = '  +0.30%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '36.48'
$ws.Range('E33').Value = '  -3.00%  '

$ws.Range('E34').Value = '  -3.94%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.88'
$ws.Range('E35').Value = '  -3.91%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0444'
$ws.Range('E36').Value = '  -1.47%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.20'
$ws.Range('E38').Value = '  -4.30%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.87'
$ws.Range('E39').Value = '  -5.14%  '

$ws.Range('E40').Value = '  -4.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.70'
$ws.Range('E41').Value = '  +0.99%  '

$ws.Range('E42').Value = '  -0.41%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.81'
$ws.Range('E43').Value = '  -3.53%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '122.62'
$ws.Range('E44').Value = '  +8.46%  '

$ws.Range('E45').Value = '  -2.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.124.38'
$ws.Range('E46').Value = '  -2.27%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.37'
$ws.Range('E47').Value = '  -4.31%  '

$ws.Range('E48').Value = '  -5.60%  '

$ws.Range('E49').Value = '  -1.43%  '

$ws.Range('E50').Value = '  -2.75%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.927'
$ws.Range('E51').Value = '  -0.94%  '
